{"js": "// Update the date line and the multiplication problems throughout the\n// document to match the \"output generated at c986bee\" refresh.\nconst replacements = [\n  [\"2024-09-15 Sunday\", \"2024-09-16 Monday\"],\n  [\"222\u00d74=\", \"212\u00d74=\"],\n  [\"123\u00d72=\", \"680\u00d78=\"],\n  [\"587\u00d79=\", \"733\u00d78=\"],\n  [\"152\u00d74=\", \"680\u00d77=\"],\n  [\"482\u00d76=\", \"853\u00d74=\"],\n  [\"536\u00d72=\", \"704\u00d72=\"],\n  [\"212\u00d76=\", \"944\u00d79=\"],\n  [\"325\u00d72=\", \"323\u00d79=\"],\n  [\"411\u00d74=\", \"500\u00d75=\"],\n  [\"246\u00d76=\", \"163\u00d78=\"],\n  [\"269\u00d72=\", \"226\u00d73=\"],\n  [\"728\u00d74=\", \"656\u00d74=\"],\n  [\"960\u00d72=\", \"252\u00d76=\"],\n  [\"157\u00d79=\", \"664\u00d72=\"],\n  [\"807\u00d79=\", \"899\u00d79=\"],\n  [\"912\u00d73=\", \"508\u00d78=\"],\n  [\"411\u00d75=\", \"180\u00d76=\"],\n  [\"529\u00d76=\", \"962\u00d72=\"],\n  [\"678\u00d79=\", \"959\u00d78=\"],\n  [\"750\u00d75=\", \"702\u00d75=\"],\n  [\"691\u00d74=\", \"236\u00d79=\"],\n  [\"740\u00d72=\", \"640\u00d74=\"],\n  [\"822\u00d78=\", \"349\u00d77=\"],\n  [\"601\u00d74=\", \"939\u00d73=\"],\n  [\"324\u00d75=\", \"777\u00d78=\"],\n];\n\nfor (const [oldText, newText] of replacements) {\n  const results = context.document.body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  for (const range of results.items) {\n    range.insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "# Update the date line and the multiplication problems throughout the\n# document to match the \"output generated at c986bee\" refresh.\n$d = $word.ActiveDocument\n\n$pairs = @(\n  @(\"2024-09-15 Sunday\", \"2024-09-16 Monday\"),\n  @(\"222\u00d74=\", \"212\u00d74=\"),\n  @(\"123\u00d72=\", \"680\u00d78=\"),\n  @(\"587\u00d79=\", \"733\u00d78=\"),\n  @(\"152\u00d74=\", \"680\u00d77=\"),\n  @(\"482\u00d76=\", \"853\u00d74=\"),\n  @(\"536\u00d72=\", \"704\u00d72=\"),\n  @(\"212\u00d76=\", \"944\u00d79=\"),\n  @(\"325\u00d72=\", \"323\u00d79=\"),\n  @(\"411\u00d74=\", \"500\u00d75=\"),\n  @(\"246\u00d76=\", \"163\u00d78=\"),\n  @(\"269\u00d72=\", \"226\u00d73=\"),\n  @(\"728\u00d74=\", \"656\u00d74=\"),\n  @(\"960\u00d72=\", \"252\u00d76=\"),\n  @(\"157\u00d79=\", \"664\u00d72=\"),\n  @(\"807\u00d79=\", \"899\u00d79=\"),\n  @(\"912\u00d73=\", \"508\u00d78=\"),\n  @(\"411\u00d75=\", \"180\u00d76=\"),\n  @(\"529\u00d76=\", \"962\u00d72=\"),\n  @(\"678\u00d79=\", \"959\u00d78=\"),\n  @(\"750\u00d75=\", \"702\u00d75=\"),\n  @(\"691\u00d74=\", \"236\u00d79=\"),\n  @(\"740\u00d72=\", \"640\u00d74=\"),\n  @(\"822\u00d78=\", \"349\u00d77=\"),\n  @(\"601\u00d74=\", \"939\u00d73=\"),\n  @(\"324\u00d75=\", \"777\u00d78=\")\n)\n\nforeach ($pair in $pairs) {\n    $old = $pair[0]\n    $new = $pair[1]\n\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $old\n    $find.Replacement.Text = $new\n    $find.Execute($old, $false, $false, $false, $false, $false, $true, 1, $false, $new, 2) | Out-Null\n}\n"}
